$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers I1 ("I0") and J1 ("IF"), matching the style of the existing
# header cells (e.g. H1) which use style index 1 (bold, centered, bordered).
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Row -> (I value, J value) for data rows 2..44
$data = @(
    @(2, 7, 8),
    @(3, 6, 7),
    @(4, 2, 3),
    @(5, 8, 8),
    @(6, 6, 7),
    @(7, 9, 9),
    @(8, 9, 9),
    @(9, 8, 9),
    @(10, 5, 7),
    @(11, 5, 7),
    @(12, 9, 9),
    @(13, 6, 6),
    @(14, 1, 5),
    @(15, 1, 5),
    @(16, 10, 10),
    @(17, 9, 9),
    @(18, 8, 9),
    @(19, 9, 9),
    @(20, 7, 8),
    @(21, 9, 9),
    @(22, 5, 5),
    @(23, 1, 3),
    @(24, 4, 4),
    @(25, 10, 11),
    @(26, 10, 10),
    @(27, 6, 6),
    @(28, 7, 8),
    @(29, 5, 6),
    @(30, 6, 7),
    @(31, 10, 10),
    @(32, 6, 6),
    @(33, 6, 6),
    @(34, 7, 8),
    @(35, 8, 8),
    @(36, 8, 8),
    @(37, 1, 2),
    @(38, 5, 5),
    @(39, 5, 6),
    @(40, 7, 8),
    @(41, 6, 7),
    @(42, 6, 7),
    @(43, 5, 6),
    @(44, 4, 5)
)

foreach ($item in $data) {
    $r = $item[0]
    $iVal = $item[1]
    $jVal = $item[2]
    $ws.Cells.Item($r, 9).Value = $iVal
    $ws.Cells.Item($r, 10).Value = $jVal
}
